$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps numeric-looking values as text,
# matching the source data which stores prices as plain text strings.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "47.425.27"
$ws.Range("E2").Value = "  +4.64%  "

$ws.Range("D3").Value = "2.498.48"
$ws.Range("E3").Value = "  +2.90%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "323.67"
$ws.Range("E5").Value = "  +1.52%  "

$ws.Range("D6").Value = "107.97"
$ws.Range("E6").Value = "  +5.10%  "

$ws.Range("D7").Value = "0.527"
$ws.Range("E7").Value = "  +2.10%  "

$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "0.544"
$ws.Range("E9").Value = "  +3.04%  "

$ws.Range("D10").Value = "38.24"
$ws.Range("E10").Value = "  +7.37%  "

$ws.Range("D11").Value = "0.0814"
$ws.Range("E11").Value = "  +1.70%  "

$ws.Range("E12").Value = "  +1.38%  "

$ws.Range("D13").Value = "18.42"
$ws.Range("E13").Value = "  +1.08%  "

$ws.Range("D14").Value = "7.21"
$ws.Range("E14").Value = "  +1.66%  "

$ws.Range("D15").Value = "2.890.65"
$ws.Range("E15").Value = "  +2.96%  "

$ws.Range("D16").Value = "2.488.03"
$ws.Range("E16").Value = "  +2.56%  "

$ws.Range("E17").Value = "  +0.67%  "

$ws.Range("D18").Value = "47.355.46"
$ws.Range("E18").Value = "  +4.76%  "

$ws.Range("D19").Value = "12.95"
$ws.Range("E19").Value = "  +5.97%  "

$ws.Range("D20").Value = "6.67"
$ws.Range("E20").Value = "  +5.10%  "

$ws.Range("E21").Value = "  +2.05%  "

$ws.Range("D22").Value = "70.72"
$ws.Range("E22").Value = "  +2.60%  "

$ws.Range("E23").Value = "  +7.30%  "

$ws.Range("D24").Value = "251.60"
$ws.Range("E24").Value = "  +2.83%  "

$ws.Range("D25").Value = "2.59"
$ws.Range("E25").Value = "  +4.18%  "

$ws.Range("D26").Value = "26.24"
$ws.Range("E26").Value = "  +1.69%  "

$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").Value = "10.08"
$ws.Range("E28").Value = "  +4.77%  "

$ws.Range("D29").Value = "35.49"
$ws.Range("E29").Value = "  +7.47%  "

$ws.Range("D30").Value = "0.139"
$ws.Range("E30").Value = "  +10.57%  "

$ws.Range("E31").Value = "  -8.43%  "

$ws.Range("E32").Value = "  +0.02%  "

$ws.Range("D33").Value = "19.85"
$ws.Range("E33").Value = "  -0.81%  "

$ws.Range("D34").Value = "5.41"
$ws.Range("E34").Value = "  +3.80%  "

$ws.Range("D35").Value = "0.0795"
$ws.Range("E35").Value = "  +3.90%  "

$ws.Range("E36").Value = "  +0.23%  "

$ws.Range("E37").Value = "  +6.41%  "

$ws.Range("D38").Value = "4.70"
$ws.Range("E38").Value = "  +5.68%  "

$ws.Range("D39").Value = "3.00"
$ws.Range("E39").Value = "  +4.04%  "

$ws.Range("E40").Value = "  +1.80%  "

$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "2.24"
$ws.Range("E41").Value = "  +2.01%  "

$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "121.72"
$ws.Range("E42").Value = "  -3.61%  "

$ws.Range("D43").Value = "21.13"
$ws.Range("E43").Value = "  +1.81%  "

$ws.Range("E44").Value = "  +2.82%  "

$ws.Range("D45").Value = "1.972.62"
$ws.Range("E45").Value = "  +1.83%  "

$ws.Range("D46").Value = "3.03"
$ws.Range("E46").Value = "  +2.86%  "

$ws.Range("E47").Value = "  -0.69%  "

$ws.Range("D48").Value = "1.81"
$ws.Range("E48").Value = "  -0.12%  "

$ws.Range("D49").Value = "9.14"
$ws.Range("E49").Value = "  +0.51%  "

$ws.Range("D50").Value = "5.29"
$ws.Range("E50").Value = "  +10.40%  "

$ws.Range("D51").Value = "79.89"
$ws.Range("E51").Value = "  +4.06%  "
